$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 5083
$ws.Range("I3").Value = 5297
$ws.Range("H4").Value = 1674
$ws.Range("I4").Value = 1215
$ws.Range("I5").Value = 492
$ws.Range("I6").Value = 5788
$ws.Range("H7").Value = 25984
$ws.Range("I7").Value = 17875

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I5").Value = 56
$ws.Range("I6").Value = 122
$ws.Range("I7").Value = 577
$ws.Range("I8").Value = 1078
$ws.Range("I14").Value = 105
$ws.Range("I16").Value = 50
$ws.Range("I19").Value = 490
$ws.Range("I22").Value = 45
$ws.Range("I23").Value = 175
$ws.Range("I29").Value = 1134
$ws.Range("I33").Value = 817
$ws.Range("I36").Value = 235
$ws.Range("I37").Value = 573
$ws.Range("I42").Value = 605
$ws.Range("I43").Value = 141
$ws.Range("I51").Value = 193
$ws.Range("I52").Value = 396
$ws.Range("I53").Value = 183
$ws.Range("I55").Value = 205
$ws.Range("I60").Value = 94
$ws.Range("H63").Value = 221
$ws.Range("I63").Value = 69
$ws.Range("I64").Value = 156
$ws.Range("I65").Value = 399
$ws.Range("I67").Value = 701
$ws.Range("I73").Value = 155
$ws.Range("I75").Value = 55
$ws.Range("I83").Value = 376
$ws.Range("I84").Value = 149
$ws.Range("I85").Value = 806
$ws.Range("I86").Value = 108
$ws.Range("I88").Value = 165
$ws.Range("I89").Value = 205
$ws.Range("I94").Value = 179
$ws.Range("I95").Value = 289
$ws.Range("I96").Value = 190
$ws.Range("I99").Value = 334
$ws.Range("H101").Value = 25984
$ws.Range("I101").Value = 17875

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 220
$ws.Range("I3").Value = 320
$ws.Range("I7").Value = 806

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I3").Value = 144
$ws.Range("I7").Value = 396

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 338
$ws.Range("I3").Value = 306
$ws.Range("I6").Value = 344
$ws.Range("I7").Value = 1078

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I6").Value = 84
$ws.Range("I7").Value = 183

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I3").Value = 179
$ws.Range("I7").Value = 577

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I3").Value = 48
$ws.Range("I7").Value = 205

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I6").Value = 71
$ws.Range("I7").Value = 190

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("I2").Value = 38
$ws.Range("I7").Value = 105

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I3").Value = 186
$ws.Range("I7").Value = 573

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I2").Value = 95
$ws.Range("I4").Value = 25
$ws.Range("I7").Value = 334

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 162
$ws.Range("I3").Value = 253
$ws.Range("I5").Value = 19
$ws.Range("I7").Value = 701

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("I3").Value = 48
$ws.Range("I6").Value = 39
$ws.Range("I7").Value = 149

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I2").Value = 132
$ws.Range("I3").Value = 118
$ws.Range("I7").Value = 399

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 133
$ws.Range("I7").Value = 376

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I3").Value = 108
$ws.Range("I7").Value = 289

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I3").Value = 306
$ws.Range("I4").Value = 36
$ws.Range("I6").Value = 254
$ws.Range("I7").Value = 817

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 339
$ws.Range("I3").Value = 392
$ws.Range("I7").Value = 1134

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I2").Value = 179
$ws.Range("I7").Value = 490

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I2").Value = 49
$ws.Range("I7").Value = 122

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 156
$ws.Range("I7").Value = 605

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I3").Value = 65
$ws.Range("I7").Value = 205

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I3").Value = 61
$ws.Range("I6").Value = 50
$ws.Range("I7").Value = 175

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("I6").Value = 54
$ws.Range("I7").Value = 156

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I5").Value = 9
$ws.Range("I7").Value = 235

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("I6").Value = 104
$ws.Range("I7").Value = 179

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I2").Value = 48
$ws.Range("I3").Value = 51
$ws.Range("I7").Value = 155

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I2").Value = 50
$ws.Range("I6").Value = 47
$ws.Range("I7").Value = 165

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("I3").Value = 15
$ws.Range("I7").Value = 56

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I2").Value = 20
$ws.Range("I4").Value = 53
$ws.Range("I7").Value = 108

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("I6").Value = 14
$ws.Range("I7").Value = 55

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I3").Value = 53
$ws.Range("I7").Value = 193

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("I2").Value = 32
$ws.Range("I7").Value = 94

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I2").Value = 29
$ws.Range("I7").Value = 141

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("I2").Value = 16
$ws.Range("I7").Value = 45

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("I2").Value = 10
$ws.Range("I7").Value = 50
